$d = $word.ActiveDocument

# 1) Merge the split "October" / " 2016" runs in the Italy row into a
#    single run reading "October 2016".
$d.Content.Find.Execute("October 2016", $true, $false, $false, $false, $false, $true, 1, $false, "October 2016", 2) | Out-Null

# 2) Add a new run "Adding new text" to the trailing paragraph, right
#    before the existing _GoBack bookmark.
$lastParagraph = $d.Paragraphs.Last
$lastParagraph.Range.InsertBefore("Adding new text")
